$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-modified date serial for each
# data row (rows 2 through 387). Bump every one of them from 45179
# (2023-09-10) to 45180 (2023-09-11).
$ws.Range("C2:C387").Value = 45180
